$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Reference cells whose style (cellXf) index we reuse via Copy+PasteSpecial
# (xlPasteFormats = -4122). These are never touched by the edits below, so
# they stay valid "donors" throughout the whole script.
#   style 1 -> header (bold, border, centered)      donor: A1
#   style 2 -> red fill   (attendance == 0)         donor: F2
#   style 3 -> yellow fill (1  <= attendance <= 19) donor: J2
#   style 4 -> white fill (20 <= attendance <= 30)  donor: L2
#   style 5 -> green fill (31 <= attendance <= 40)  donor: H21
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

function Copy-Style($donorA1, $targetA1) {
    $ws.Range($donorA1).Copy() | Out-Null
    $ws.Range($targetA1).PasteSpecial($xlPasteFormats) | Out-Null
}

function Style-ForAttendance($n) {
    if ($n -eq 0) { return "F2" }
    elseif ($n -ge 1 -and $n -le 19) { return "J2" }
    elseif ($n -ge 20 -and $n -le 30) { return "L2" }
    else { return "H21" }
}

# ---------------------------------------------------------------------------
# 1. New header cells X1 ("05-02_A") / Y1 ("05-02_0"), styled like the rest
#    of row 1 (style 1).
# ---------------------------------------------------------------------------
Copy-Style "A1" "X1"
Copy-Style "A1" "Y1"
$ws.Range("X1").Value = "05-02_A"
$ws.Range("Y1").Value = "05-02_0"

# ---------------------------------------------------------------------------
# 2. Per-row data. Columns:
#    row | isEmpty | newWValue | wAlreadyNumeric | xValue | yValue | eOverride
#    - isEmpty=1          -> no 05-02 game for that player: X gets the
#                            "white/no-data" style (4) with no value, Y stays
#                            blank.
#    - newWValue          -> the W column (05-01_0 total) used to be stored
#                            as text; it is rewritten as a genuine number
#                            with the same value (skipped when the cell is
#                            already numeric).
#    - wAlreadyNumeric=1  -> W already holds a numeric value; leave it alone.
#    - xValue / yValue    -> the new 05-02_A attendance count / 05-02_0
#                            running total.
#    - eOverride          -> (row 76 only) the "Now" library column changes
#                            from 二馆 to 一馆.
# ---------------------------------------------------------------------------
$rowData = @"
2|0|3083|0|0|3083|
3|0|0|0|0|0|
4|0|0|0|0|0|
5|0|2506|0|0|2499|
6|0|2720|0|0|2772|
7|1||0|||
8|1||0|||
9|1||0|||
10|1||0|||
11|1||0|||
12|1||0|||
13|1||0|||
14|1||0|||
15|1||0|||
16|1||0|||
17|1||0|||
18|0|3995|0|11|4066|
19|0|2636|0|0|2652|
20|0|4213|0|27|4365|
21|0|4708|0|40|4912|
22|0|4681|0|20|4863|
23|0|5283|0|23|5537|
24|0|4832|0|33|4974|
25|0|5118|0|30|5288|
26|0|2778|0|6|2934|
27|0|2500|0|0|2500|
28|1||0|||
29|0|3050|0|21|3623|
30|0|4956|0|16|5010|
31|0|4862|0|30|4990|
32|0|2696|0|0|2696|
33|1||0|||
34|0|2500|0|0|2500|
35|0|4601|0|20|4840|
36|0|2727|0|0|2722|
37|0|4541|0|30|4731|
38|0|4918|0|22|5177|
39|0|4575|0|25|4683|
40|0|0|0|0|0|
41|0|4281|0|20|4432|
42|0|3044|0|0|3149|
43|1||0|||
44|1||1|||
45|0|3987|0|28|4083|
46|0|3915|0|15|4061|
47|0|5226|0|33|5398|
48|0|4817|0|30|5246|
49|0|4695|0|30|4743|
50|0|4777|0|23|4982|
51|0|3875|0|0|3922|
52|0|4949|0|30|5055|
53|0|3635|0|5|3751|
54|0|4623|0|20|4733|
55|0|3705|0|0|3701|
56|0|5165|0|30|5232|
57|0|4203|0|16|4201|
58|0|4106|0|23|4186|
59|0|4074|0|20|4131|
60|0|4206|0|20|4262|
61|1||0|||
62|0|3988|0|30|4027|
63|0|4035|0|20|4105|
64|0|4143|0|0|4133|
65|0|3888|0|23|3995|
66|0|0|0|0|0|
67|0|0|0|0|0|
68|0|2502|0|0|2499|
69|0|2956|0|0|2963|
70|0|1518|0|0|1510|
71|0|0|0|0|0|
72|0|0|0|0|0|
73|0|2605|0|0|2628|
74|0|0|0|0|0|
75|0|0|0|0|0|
76|0||1|0|4510|一馆
77|0|3675|0|20|3929|
78|0|2862|0|0|2909|
79|0|1306|0|0|1305|
80|0|0|0|0|0|
81|0|2647|0|0|2642|
82|0|0|0|0|0|
83|0|0|0|0|0|
84|0|1524|0|0|1515|
85|0|0|0|0|0|
86|0|0|0|0|0|
87|0|0|0|0|0|
88|0|0|0|0|0|
89|0|0|0|0|0|
90|0|0|0|0|0|
91|0|0|0|0|0|
92|0|0|0|0|0|
93|0|0|0|0|0|
94|0|0|0|0|0|
95|0|0|0|0|0|
96|0|2790|0|0|2785|
97|0|0|0|0|0|
98|0|0|0|0|0|
99|0|0|0|0|0|
100|0|0|0|0|0|
101|0|0|0|0|0|
102|0|0|0|0|0|
103|0|0|0|0|0|
104|0|0|0|0|0|
105|0|0|0|0|0|
106|1||0|||
107|1||0|||
108|1||0|||
109|1||0|||
110|1||0|||
111|1||0|||
112|1||0|||
113|1||0|||
114|1||0|||
115|1||1|||
116|1||1|||
117|1||1|||
118|0|5403|0|20|5578|
119|0|1609|0|0|1641|
"@

$lines = $rowData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split '\|'
    $r          = [int]$parts[0]
    $isEmpty    = [int]$parts[1]
    $newW       = $parts[2]
    $wAlreadyN  = [int]$parts[3]
    $xVal       = $parts[4]
    $yVal       = $parts[5]
    $eOverride  = $parts[6]

    if ($eOverride -ne "") {
        $ws.Range("E$r").Value = $eOverride
    }

    if ($wAlreadyN -eq 0) {
        # W currently holds its number as text -> rewrite as a real number
        $ws.Range("W$r").Value = [double]$newW
    }

    if ($isEmpty -eq 1) {
        Copy-Style "L2" "X$r"
    } else {
        $n = [int]$xVal
        $donor = Style-ForAttendance $n
        Copy-Style $donor "X$r"
        $ws.Range("X$r").Value = $n
        $ws.Range("Y$r").Value = "'" + $yVal
        $ws.Range("Y$r").ClearFormats()
    }
}

# ---------------------------------------------------------------------------
# 3. Brand-new row 120 (a player who only appears starting 05-02).
# ---------------------------------------------------------------------------
$ws.Range("B120").Value = "Hong"
$ws.Range("E120").Value = "三馆"
foreach ($col in @("F","H","J","L","N","P","R","T","V")) {
    Copy-Style "L2" "$col`120"
}
Copy-Style "J2" "X120"
$ws.Range("X120").Value = 18
$ws.Range("Y120").Value = "'1851"
$ws.Range("Y120").ClearFormats()

# A120 must stay literal text "59304163" (not a number) with the default
# (unstyled) cell format, matching the source data's inlineStr typing.
$ws.Range("A120").Value = "'59304163"
$ws.Range("A120").ClearFormats()
